$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-25 Monday" "2023-12-26 Tuesday"

Replace-Text "86÷9=" "49÷4="
Replace-Text "41÷9=" "29÷5="
Replace-Text "58÷4=" "25÷6="
Replace-Text "62÷5=" "22÷8="
Replace-Text "55÷7=" "92÷5="
Replace-Text "75÷3=" "97÷4="
Replace-Text "70÷7=" "98÷8="
Replace-Text "91÷3=" "87÷6="
Replace-Text "92÷9=" "58÷5="
Replace-Text "99÷7=" "71÷3="
Replace-Text "65÷2=" "80÷6="
Replace-Text "63÷8=" "66÷2="
Replace-Text "60÷2=" "77÷4="
Replace-Text "91÷7=" "24÷8="
Replace-Text "55÷2=" "26÷7="
Replace-Text "73÷7=" "36÷4="
Replace-Text "40÷8=" "93÷5="
Replace-Text "15÷3=" "25÷4="
Replace-Text "21÷2=" "89÷4="
Replace-Text "54÷5=" "37÷9="
Replace-Text "68÷6=" "58÷8="
Replace-Text "93÷7=" "45÷2="
Replace-Text "75÷7=" "59÷5="
Replace-Text "30÷4=" "72÷2="
Replace-Text "18÷9=" "50÷4="
